$wb = $excel.ActiveWorkbook

$wsSingle = $wb.Worksheets.Item("single")
$wsMulti = $wb.Worksheets.Item("multi")

# Delete column G (Other_indexes) on both "single" and "multi" sheets
$wsSingle.Columns.Item("G").Delete()
$wsMulti.Columns.Item("G").Delete()

# Fix the "single" sheet formula that incorrectly said "_Multi"
$wsSingle.Range("B7").Formula = '=VLOOKUP(B$5, config!$B$4:$E$14,2,FALSE) & "_Single"'
$wsSingle.Range("K7").Formula = '=VLOOKUP(K$5, config!$B$4:$E$14,2,FALSE) & " - Single"'
